$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.767.18'
$ws.Range("E2").Value = '  +3.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.469.88'
$ws.Range("E3").Value = '  +4.21%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.61'
$ws.Range("E5").Value = '  +4.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.72'
$ws.Range("E6").Value = '  +4.00%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.468.59'
$ws.Range("E8").Value = '  +4.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.566'
$ws.Range("E9").Value = '  +7.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.58'
$ws.Range("E10").Value = '  +1.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.126'
$ws.Range("E11").Value = '  +6.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.452'
$ws.Range("E12").Value = '  +3.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.066.69'
$ws.Range("E13").Value = '  +4.25%  '
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000199'
$ws.Range("E15").Value = '  +10.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.04'
$ws.Range("E16").Value = '  +3.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.757.77'
$ws.Range("E17").Value = '  +3.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.467.65'
$ws.Range("E18").Value = '  +4.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.47'
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.45'
$ws.Range("E20").Value = '  +4.52%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '398.91'
$ws.Range("E21").Value = '  +3.35%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.59'
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.550'
$ws.Range("E23").Value = '  +2.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.32'
$ws.Range("E24").Value = '  +3.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000123'
$ws.Range("E26").Value = '  +26.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.65'
$ws.Range("E27").Value = '  +9.19%  '
$ws.Range("E28").Value = '  +1.17%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.05'
$ws.Range("E30").Value = '  +8.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.77'
$ws.Range("E31").Value = '  +4.71%  '
$ws.Range("E32").Value = '  +3.37%  '
$ws.Range("E33").Value = '  +5.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.95'
$ws.Range("E34").Value = '  +4.26%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.04'
$ws.Range("E36").Value = '  +4.39%  '
$ws.Range("E37").Value = '  -0.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.69'
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0789'
$ws.Range("E39").Value = '  +7.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.89'
$ws.Range("E40").Value = '  -0.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '27.77'
$ws.Range("E41").Value = '  +2.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.913.47'
$ws.Range("E42").Value = '  +1.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0326'
$ws.Range("E43").Value = '  +3.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.47'
$ws.Range("E44").Value = '  +2.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.35'
$ws.Range("E45").Value = '  +4.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.774'
$ws.Range("E46").Value = '  +2.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.89'
$ws.Range("E47").Value = '  +8.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.10'
$ws.Range("E48").Value = '  +5.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.22'
$ws.Range("E49").Value = '  +25.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.60'
$ws.Range("E50").Value = '  +4.68%  '
$ws.Range("E51").Value = '  +2.61%  '
